$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new daily rows (245-247) to the Concordia report, reusing the
# existing formatting of column A (date style) from the last current row.

$rows = @(
    @{ Row = 245; A = 44319; B = 2; C = 8; D = 97.22897423432183 },
    @{ Row = 246; A = 44320; B = 0; C = 8; D = 97.22897423432183 },
    @{ Row = 247; A = 44321; B = 0; C = 6; D = 72.92173067574137 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    # Copy column A's formatting from the row above so the new date cell
    # keeps the same style (centered, bordered, date number format).
    $ws.Range("A" + ($rowNum - 1)).Copy()
    $ws.Range("A" + $rowNum).PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    $ws.Cells.Item($rowNum, 1).Value = $r.A
    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D
}
